# Logged 2021 divisional round, simulated season from conference round
$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 (Home team) updated with new simulated totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 472
$wsOff.Range("C2").Value = 337
$wsOff.Range("D2").Value = 97
$wsOff.Range("E2").Value = 45
$wsOff.Range("F2").Value = 8
$wsOff.Range("G2").Value = 9

# --- DEF sheet: row 2 (Home team) updated with new simulated totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 511
$wsDef.Range("C2").Value = 344
$wsDef.Range("D2").Value = 122
$wsDef.Range("E2").Value = 62
$wsDef.Range("F2").Value = 11
$wsDef.Range("G2").Value = 7
